$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4632.5
$ws.Range("I19").Value = 4131.75
$ws.Range("K19").Value = 4131.75
$ws.Range("M19").Value = -3956.75
$ws.Range("H38").Value = 11381.5
$ws.Range("I38").Value = 18765
$ws.Range("J38").Value = 3998
$ws.Range("K38").Value = 56295
$ws.Range("L38").Value = 11994
$ws.Range("M38").Value = -55923
$ws.Range("N38").Value = -12738
$ws.Range("H39").Value = 2384.3333
$ws.Range("I39").Value = 77.5
$ws.Range("K39").Value = 232.5
$ws.Range("M39").Value = 63.5
$ws.Range("H40").Value = 35374.875
$ws.Range("I40").Value = 57499.5
$ws.Range("J40").Value = 28000
$ws.Range("K40").Value = 57499.5
$ws.Range("L40").Value = 28000
$ws.Range("M40").Value = -57324.5
$ws.Range("N40").Value = -28350
$ws.Range("H58").Value = 1620.8334
$ws.Range("I58").Value = 45
$ws.Range("J58").Value = 9500
$ws.Range("K58").Value = 135
$ws.Range("L58").Value = 28500
$ws.Range("M58").Value = 15
$ws.Range("N58").Value = -28800
$ws.Range("H86").Value = 73733864
$ws.Range("I86").Value = 93842780
$ws.Range("J86").Value = 1166.6666
$ws.Range("K86").Value = 93842780
$ws.Range("L86").Value = 1166.6666
$ws.Range("M86").Value = -93841657
$ws.Range("N86").Value = -3412.6666
$ws.Range("H89").Value = 73733864
$ws.Range("I89").Value = 93842780
$ws.Range("J89").Value = 1166.6666
$ws.Range("K89").Value = 469213900
$ws.Range("L89").Value = 5833.333000000001
$ws.Range("M89").Value = -469208284
$ws.Range("N89").Value = -17065.333
$ws.Range("H137").Value = 3706.9656
$ws.Range("I137").Value = 2737.3076
$ws.Range("K137").Value = 8211.9228
$ws.Range("M137").Value = -5661.9228
$ws.Range("H138").Value = 5601.3057
$ws.Range("I138").Value = 1928
$ws.Range("J138").Value = 5996.892
$ws.Range("K138").Value = 5784
$ws.Range("L138").Value = 17990.676
$ws.Range("M138").Value = -644
$ws.Range("N138").Value = -28270.676
$ws.Range("H140").Value = 68038.78
$ws.Range("J140").Value = 67793.625
$ws.Range("L140").Value = 67793.625
$ws.Range("N140").Value = -78153.625

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3461.0876
$ws.Range("I32").Value = 1865.4348
$ws.Range("K32").Value = 1865.4348
$ws.Range("M32").Value = -1578.4348
$ws.Range("H74").Value = 3742.075
$ws.Range("I74").Value = 1507.625
$ws.Range("K74").Value = 1507.625
$ws.Range("M74").Value = -633.625
$ws.Range("H77").Value = 3742.075
$ws.Range("I77").Value = 1507.625
$ws.Range("K77").Value = 7538.125
$ws.Range("M77").Value = -3170.125

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 166674340
$ws.Range("I105").Value = 333342660
$ws.Range("K105").Value = 333342660
$ws.Range("M105").Value = -333340913

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1487.2858
$ws.Range("I16").Value = 1551.8334
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 1551.8334
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -1264.8334
$ws.Range("N16").Value = -1674
$ws.Range("H107").Value = 677643.7
$ws.Range("I107").Value = 1012979.4
$ws.Range("K107").Value = 1012979.4
$ws.Range("M107").Value = -1011059.4
$ws.Range("H113").Value = 1487.2858
$ws.Range("I113").Value = 1551.8334
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 1551.8334
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 618.1666
$ws.Range("N113").Value = -5440

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 291520.44
$ws.Range("I68").Value = 1374.2
$ws.Range("J68").Value = 462194.72
$ws.Range("K68").Value = 4122.6
$ws.Range("L68").Value = 1386584.16
$ws.Range("M68").Value = -3311.6
$ws.Range("N68").Value = -1388206.16
$ws.Range("H71").Value = 291520.44
$ws.Range("I71").Value = 1374.2
$ws.Range("J71").Value = 462194.72
$ws.Range("K71").Value = 12367.8
$ws.Range("L71").Value = 4159752.48
$ws.Range("M71").Value = -8311.800000000001
$ws.Range("N71").Value = -4167864.48
$ws.Range("H107").Value = 431340.9
$ws.Range("I107").Value = 1377.6666
$ws.Range("J107").Value = 569543.4
$ws.Range("K107").Value = 4132.9998
$ws.Range("L107").Value = 1708630.2
$ws.Range("M107").Value = -2212.9998
$ws.Range("N107").Value = -1712470.2
$ws.Range("H131").Value = 16284297
$ws.Range("J131").Value = 7445182
$ws.Range("L131").Value = 22335546
$ws.Range("N131").Value = -22345626
$ws.Range("H132").Value = 2266.111
$ws.Range("I132").Value = 933.3333
$ws.Range("J132").Value = 2932.5
$ws.Range("K132").Value = 8399.9997
$ws.Range("L132").Value = 26392.5
$ws.Range("M132").Value = -5869.9997
$ws.Range("N132").Value = -31452.5
$ws.Range("H134").Value = 3421.0527
$ws.Range("I134").Value = 2937.5
$ws.Range("K134").Value = 8812.5
$ws.Range("M134").Value = -3742.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 309865.03
$ws.Range("I122").Value = 614902.0600000001
$ws.Range("K122").Value = 1844706.18
$ws.Range("M122").Value = -1842256.18

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4039.7576
$ws.Range("I46").Value = 1687
$ws.Range("J46").Value = 4562.593
$ws.Range("K46").Value = 1687
$ws.Range("L46").Value = 4562.593
$ws.Range("M46").Value = -1499
$ws.Range("N46").Value = -4938.593
$ws.Range("H118").Value = 15000
$ws.Range("J118").Value = 15000
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = -18314
$ws.Range("H122").Value = 4504.4165
$ws.Range("I122").Value = 3604.8
$ws.Range("J122").Value = 9002.5
$ws.Range("K122").Value = 10814.4
$ws.Range("L122").Value = 27007.5
$ws.Range("M122").Value = -8364.400000000001
$ws.Range("N122").Value = -31907.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13780.333
$ws.Range("J41").Value = 13780.333
$ws.Range("L41").Value = 13780.333
$ws.Range("N41").Value = -14560.333
$ws.Range("H81").Value = 2981489.8
$ws.Range("I81").Value = 4172085.8
$ws.Range("J81").Value = 4999.5
$ws.Range("K81").Value = 8344171.6
$ws.Range("L81").Value = 9999
$ws.Range("M81").Value = -8343110.6
$ws.Range("N81").Value = -12121
$ws.Range("H84").Value = 2981489.8
$ws.Range("I84").Value = 4172085.8
$ws.Range("J84").Value = 4999.5
$ws.Range("K84").Value = 41720858
$ws.Range("L84").Value = 49995
$ws.Range("M84").Value = -41715554
$ws.Range("N84").Value = -60603
$ws.Range("H126").Value = 2152.4546
$ws.Range("I126").Value = 1452.6471
$ws.Range("K126").Value = 4357.9413
$ws.Range("M126").Value = -1887.9413
$ws.Range("H132").Value = 1161665.9
$ws.Range("I132").Value = 1393217.9
$ws.Range("K132").Value = 4179653.7
$ws.Range("M132").Value = -4177123.7
